$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C24 / C32 were previously blank; give them "Image3" and pick up the
# highlighted style already used by the neighbouring C23/C31 cells
# (same fill/font as style index 6 in the sheet: orange fill).
# These are written before the "Image2" cells below so the two brand
# new shared-string entries land in the authored order: Image3, Image2.
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("C24").Value = "Image3"

$ws.Range("C31").Copy() | Out-Null
$ws.Range("C32").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("C32").Value = "Image3"

# Update shared-string-backed values in column C for the two existing
# "Image1" cells, which become "Image2"
$ws.Range("C23").Value = "Image2"
$ws.Range("C31").Value = "Image2"

# Update the selection recorded with the sheet view
$ws.Activate()
$ws.Range("B15").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
